$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw input values in row 6 (L3 leg)
$ws.Range("C6").Value = 1180
$ws.Range("F6").Value = 1700
$ws.Range("G6").Value = 2150

# Update the selection to match the new active cell
$ws.Range("H15").Select()
